$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.341.49"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "1.844.25"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.05"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6274"
$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07460"
$ws.Range("E8").Value = "  -1.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2896"
$ws.Range("E9").Value = "  -0.46%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.48"
$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "1.844.83"
$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6790"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("E15").Value = "  -0.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.87"
$ws.Range("E16").Value = "  -1.45%  "

$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").Value = "29.376.39"
$ws.Range("E18").Value = "  +0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.18"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.498"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.81"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.449"
$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1365"
$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.50"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06482"
$ws.Range("E28").Value = "  +15.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.417"
$ws.Range("E29").Value = "  -0.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.484"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.087"
$ws.Range("E31").Value = "  -0.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.083"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("E33").Value = "  -0.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.141"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6942"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.583"
$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "1.261.11"
$ws.Range("E37").Value = "  +2.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.833"
$ws.Range("E38").Value = "  +3.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01832"
$ws.Range("E39").Value = "  +1.92%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.789"
$ws.Range("E40").Value = "  +6.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9188"
$ws.Range("E41").Value = "  +1.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9989"
$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("D43").Value = "2.002.66"
$ws.Range("E43").Value = "  +1.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.57"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.03"
$ws.Range("E45").Value = "  +0.81%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.736"
$ws.Range("E46").Value = "  +3.00%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.060"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000117"
$ws.Range("E48").Value = "  +0.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1157"
$ws.Range("E49").Value = "  +0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.996"
$ws.Range("E50").Value = "  -0.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3943"
$ws.Range("E51").Value = "  -1.20%  "
